# "change check_combination to line": the "criteria" column (E) is
# renamed to "minor" and four new single-flag columns are appended
# (ML/DS, internship, no mention degree, Bootcamp), each populated
# with a 1 on the rows that match. Row heights shrink off the
# 409.6-clamp to their natural autofit heights, the header row is
# frozen, the active selection moves to G25, and a handful of column
# widths / the page setup are refreshed to match re-saving in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: rename "criteria" -> "minor", add 4 new headers ---
$ws.Range("E1").Value = "minor"
$ws.Range("F1").Value = "ML/DS"
$ws.Range("G1").Value = "internship"
$ws.Range("H1").Value = "no mention degree"
$ws.Range("I1").Value = "Bootcamp"

# --- new per-row flag values (1 = applies) ---
$ws.Range("G4").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("I25").Value = 1
$ws.Range("E26").Value = 1

# --- row heights: drop from the old 409.6 clamp to natural heights ---
$ws.Rows.Item(2).RowHeight = 171
$ws.Rows.Item(3).RowHeight = 171
$ws.Rows.Item(4).RowHeight = 304
$ws.Rows.Item(6).RowHeight = 304
$ws.Rows.Item(7).RowHeight = 323
$ws.Rows.Item(9).RowHeight = 380
$ws.Rows.Item(10).RowHeight = 266
$ws.Rows.Item(11).RowHeight = 114
$ws.Rows.Item(12).RowHeight = 152
$ws.Rows.Item(13).RowHeight = 171
$ws.Rows.Item(17).RowHeight = 285
$ws.Rows.Item(19).RowHeight = 171
$ws.Rows.Item(21).RowHeight = 323
$ws.Rows.Item(22).RowHeight = 95
$ws.Rows.Item(25).RowHeight = 304
$ws.Rows.Item(26).RowHeight = 285
$ws.Rows.Item(29).RowHeight = 247

# --- column widths: narrow old col E, add widths for new cols F & H ---
$ws.Columns.Item(5).ColumnWidth = 6.71058125
$ws.Columns.Item(6).ColumnWidth = 7.425425
$ws.Columns.Item(8).ColumnWidth = 18.9957375

# --- freeze the header row and move the selection down to the new data ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G25").Select()

# --- page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
